# Add new columns I ("I0") and J ("IF") to the sheet, mirroring the header
# style used by the existing H ("IP") column, and populate data rows 2-44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header cell format (bold font, border, centered/top alignment)
# from H1 into I1:J1, then overwrite with the new header text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..44 (columns I and J), taken from the target diff.
$iVals = @(9,9,9,8,9,9,9,9,10,9,7,8,8,6,5,9,8,5,5,7,8,9,7,8,9,9,9,8,8,7,8,9,7,8,9,9,8,9,7,9,8,8,9)
$jVals = @(10,9,9,8,9,9,9,9,10,9,7,8,8,6,5,9,8,5,6,8,8,9,8,8,9,9,9,8,8,7,9,9,8,8,9,9,8,9,7,9,8,8,9)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
